$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = "Kin Id"
$ws.Range("A3").Value = "61652_FS"

$ws.Range("A3").Select()
